# "Generate Report for Archive"
#
# 1) Status text "Ready for handoff" -> "In Translation" everywhere it
#    appears (Overview sheet E2/F2 per-locale status roll-up, and the
#    "Status" column on each locale sheet).
# 2) Narrow the "Status" column(s) from ~17.22 chars down to ~13.41 chars
#    on the Overview sheet (columns E & F) and on each locale sheet
#    (column C).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- 1) Update the status value wherever it is used ---
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- 2) Narrow the Status columns ---
# Target stored column width ~= 13.4101845877511 chars (was ~17.2159881591797).
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
